$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja4")

$ws.Cells.Item(41, 8).Value = 0.35
$ws.Cells.Item(41, 10).Value = 17.26027397260274
$ws.Cells.Item(41, 11).Value = 119.6054794520548
$ws.Cells.Item(41, 12).Value = 0.3975555555555556

$ws.Cells.Item(42, 8).Value = 0.355
$ws.Cells.Item(42, 10).Value = 20.424657534246574
$ws.Cells.Item(42, 11).Value = 122.83315068493151
$ws.Cells.Item(42, 12).Value = 0.3968619047619047

$ws.Cells.Item(43, 8).Value = 0.36
$ws.Cells.Item(43, 10).Value = 26.630136986301366
$ws.Cells.Item(43, 11).Value = 129.16273972602738
$ws.Cells.Item(43, 12).Value = 0.3942370370370369

$ws.Cells.Item(53, 8).Value = 0.1925
$ws.Cells.Item(53, 10).Value = 9.493150684931507
$ws.Cells.Item(53, 11).Value = 121.42790410958905
$ws.Cells.Item(53, 12).Value = 0.4345102777777781

$ws.Cells.Item(54, 8).Value = 0.275
$ws.Cells.Item(54, 10).Value = 20.34246575342466
$ws.Cells.Item(54, 11).Value = 133.45979452054797
$ws.Cells.Item(54, 12).Value = 0.4523268518518522

$ws.Cells.Item(55, 8).Value = 0.3275
$ws.Cells.Item(55, 10).Value = 32.3013698630137
$ws.Cells.Item(55, 11).Value = 146.72221917808218
$ws.Cells.Item(55, 12).Value = 0.4737113888888889

$ws.Cells.Item(56, 8).Value = 0.1925
$ws.Cells.Item(56, 10).Value = 9.493150684931507
$ws.Cells.Item(56, 11).Value = 121.42790410958905
$ws.Cells.Item(56, 12).Value = 0.4345102777777781

$ws.Cells.Item(57, 8).Value = 0.275
$ws.Cells.Item(57, 10).Value = 20.34246575342466
$ws.Cells.Item(57, 11).Value = 133.45979452054797
$ws.Cells.Item(57, 12).Value = 0.4523268518518522

$ws.Cells.Item(58, 8).Value = 0.3275
$ws.Cells.Item(58, 10).Value = 32.3013698630137
$ws.Cells.Item(58, 11).Value = 146.72221917808218
$ws.Cells.Item(58, 12).Value = 0.4737113888888889

$ws.Cells.Item(102, 8).Value = 0.4123148148148158
$ws.Cells.Item(102, 10).Value = 6.777777777777794
$ws.Cells.Item(102, 11).Value = 108.91333333333336
$ws.Cells.Item(102, 12).Value = 0.5422277777777789

$ws.Cells.Item(103, 8).Value = 0.4192242798353917
$ws.Cells.Item(103, 10).Value = 10.337037037037057
$ws.Cells.Item(103, 11).Value = 112.54377777777779
$ws.Cells.Item(103, 12).Value = 0.5087198765432107

$ws.Cells.Item(104, 8).Value = 0.44089196437281025
$ws.Cells.Item(104, 10).Value = 21.742617421124887
$ws.Cells.Item(104, 11).Value = 124.17746976954739
$ws.Cells.Item(104, 12).Value = 0.4902653592158218

$ws.Cells.Item(105, 8).Value = 0.5310255745122918
$ws.Cells.Item(105, 10).Value = 39.28134386803255
$ws.Cells.Item(105, 11).Value = 142.0669707453932
$ws.Cells.Item(105, 12).Value = 0.5686831230395748

$ws.Cells.Item(106, 8).Value = 0.5631689748637602
$ws.Cells.Item(106, 10).Value = 55.54543313724758
$ws.Cells.Item(106, 11).Value = 158.65634179999253
$ws.Cells.Item(106, 12).Value = 0.5947101321388131

$ws.Cells.Item(107, 8).Value = 0.38
$ws.Cells.Item(107, 10).Value = 6.246575342465753
$ws.Cells.Item(107, 11).Value = 108.37150684931507
$ws.Cells.Item(107, 12).Value = 0.5092666666666668

$ws.Cells.Item(108, 8).Value = 0.38
$ws.Cells.Item(108, 10).Value = 9.36986301369863
$ws.Cells.Item(108, 11).Value = 111.55726027397259
$ws.Cells.Item(108, 12).Value = 0.46871111111111063

$ws.Cells.Item(109, 8).Value = 0.38
$ws.Cells.Item(109, 10).Value = 18.73972602739726
$ws.Cells.Item(109, 11).Value = 121.1145205479452
$ws.Cells.Item(109, 12).Value = 0.42815555555555546

$ws.Cells.Item(110, 8).Value = 0.38
$ws.Cells.Item(110, 10).Value = 28.10958904109589
$ws.Cells.Item(110, 11).Value = 130.6717808219178
$ws.Cells.Item(110, 12).Value = 0.4146370370370368

$ws.Cells.Item(111, 8).Value = 0.38
$ws.Cells.Item(111, 10).Value = 37.47945205479452
$ws.Cells.Item(111, 11).Value = 140.2290410958904
$ws.Cells.Item(111, 12).Value = 0.40787777777777773

$ws.Cells.Item(112, 8).Value = 0.3
$ws.Cells.Item(112, 10).Value = 4.931506849315069
$ws.Cells.Item(112, 11).Value = 107.03013698630137
$ws.Cells.Item(112, 12).Value = 0.4276666666666668

$ws.Cells.Item(113, 8).Value = 0.3
$ws.Cells.Item(113, 10).Value = 7.397260273972603
$ws.Cells.Item(113, 11).Value = 109.54520547945205
$ws.Cells.Item(113, 12).Value = 0.3871111111111112

$ws.Cells.Item(114, 8).Value = 0.3
$ws.Cells.Item(114, 10).Value = 14.794520547945206
$ws.Cells.Item(114, 11).Value = 117.09041095890412
$ws.Cells.Item(114, 12).Value = 0.3465555555555556

$ws.Cells.Item(115, 8).Value = 0.3
$ws.Cells.Item(115, 10).Value = 22.191780821917806
$ws.Cells.Item(115, 11).Value = 124.63561643835618
$ws.Cells.Item(115, 12).Value = 0.33303703703703735

$ws.Cells.Item(116, 8).Value = 0.3
$ws.Cells.Item(116, 10).Value = 29.589041095890412
$ws.Cells.Item(116, 11).Value = 132.18082191780823
$ws.Cells.Item(116, 12).Value = 0.3262777777777778
